$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.498.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.68%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.620.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.67%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'535.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.78%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'142.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.87%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.568"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.41%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'6.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +5.50%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -2.19%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.334"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.17%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +1.17%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'3.082.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.62%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'58.419.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.67%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'20.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.80%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'2.619.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.22%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  -1.46%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'4.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.72%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'333.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.21%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'10.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.28%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  -1.90%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +0.16%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'66.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.55%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  +1.43%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("B25").Value = "'Kaspa"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'0.163"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.49%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("B26").Value = "'Binance-PegBSC-USD"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.09%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'7.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.59%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.01%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'0.0₃0733"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.73%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'1.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.09%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'5.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.22%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'18.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.11%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'150.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.53%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'3.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.44%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.847"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.73%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -2.15%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.812"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.62%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  -3.15%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'3.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.69%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'281.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.94%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.01%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.593"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.23%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'10.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.52%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'18.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.95%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.0526"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.50%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.0934"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.85%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  +0.01%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.938.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.02%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'4.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.42%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'17.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.64%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'113.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.40%  "
$ws.Range("E51").Style = "Normal"
